$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 data: description label and numeric date value
$ws.Range("A3").Value = "TAZ_V4_date"
$ws.Range("B3").Value = 240404

# Move the active cell selection from B5 to B4
$ws.Range("B4").Select()
